# "Generate Report for Archive" - refresh the localization-status report:
#   1. The single localization job that was "Ready for handoff" has moved on
#      to "In Translation", so update its Status everywhere it's shown
#      (Overview!E2/F2 per-language status columns, and the Status column
#      on each per-language detail sheet).
#   2. Re-fit the now-narrower Status/zh-cn/de-de columns so the report
#      doesn't carry the old, too-wide column width around.

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "In Translation"
$overview.Range("F2").Value = "In Translation"
$overview.Columns.Item(5).ColumnWidth = 12.5
$overview.Columns.Item(6).ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = "In Translation"
$zhcn.Columns.Item(3).ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = "In Translation"
$dede.Columns.Item(3).ColumnWidth = 12.5
